# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2410   (columns A-J)
#   *_new -> *_FV2504   (columns L-U)
# then turn the used range into an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $colOld = [char]([int][char]'A' + $i)   # A..J
    $colNew = [char]([int][char]'L' + $i)   # L..U
    $ws.Range("$($colOld)1").Value = "$($labels[$i])_FV2410"
    $ws.Range("$($colNew)1").Value = "$($labels[$i])_FV2504"
}

# Turn A1:U82 into a native Excel Table ("Table1") with the header row.
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U82"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
